$wb = $excel.ActiveWorkbook

# Sheet "609400404223983664": C1 id update
$ws1 = $wb.Worksheets.Item("609400404223983664")
$ws1.Range("C1").NumberFormat = "@"
$ws1.Range("C1").Value = "1135594725156601937"

# Sheet "609400404223983664_youtube": C1 id update + append new video ids
$ws2 = $wb.Worksheets.Item("609400404223983664_youtube")
$ws2.Range("C1").NumberFormat = "@"
$ws2.Range("C1").Value = "1135594725156601937"
$ws2.Range("G2").Value = "ga2BNK30Pv0,Hh-Gny6H8fw,Sw3HToXF2uc,4-nLzR9GVu0,ZJAs69NY__k,61pCygus9Ok,D7aLMfQwqCY,-jZW07fxQBs,ovkJ7l7dCKI,5ywARf7L35U"
$ws2.Range("C5").Value = "None,kLYOmpR-WTs,YGXF14e3XS4"
$ws2.Range("G5").Value = "mvLCbwzf7Fw,SizDUmbnWNk,XA8SQ49eULc,2fbuxe_K0NI,OJMHrnZUMTE,wjCANhQFs2k,J8dvBe4rGz0,jcynQNgVGlA,T-gfBgxFZX8,GraWAA-5lGA,7tmWUCDgEJk,XU46ItenHuI,HEE4gtFvA-E,oCFvnaYW0qI,pHYMlt8j4yA,0Y5fYqlkclI,LMW9zpzQdEw,VSy_q-fi7i0,_TSV8mI0PFw,LXv4w2q3SnI,FKBBJbwQOms,Wza0o6ptS-A,YGXF14e3XS4,LLecLw2kSlQ,QWOMNNR4Z5g,HJMT2vA3WYo,GAAhG5c86RA,lBvNe5YF-yE,sBF0UbN3TxQ"

# Sheet "987051677045624912_youtube": append new video id
$ws8 = $wb.Worksheets.Item("987051677045624912_youtube")
$ws8.Range("G2").Value = "vppYemSrDfk,iQ5u1GMOGc8,5MMWKgVbEgk,3Q2LZ-YJhp0,6n0hbcD0ooU"
